# Quarterly indexing esoteric bug-fix: shift each date in column A (rows 2-73)
# from the 1st of its month to the 15th of the FOLLOWING month.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2

    $d = [DateTime]::FromOADate($serial)
    $d2 = $d.AddMonths(1)
    $d3 = Get-Date -Year $d2.Year -Month $d2.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $newSerial = [Math]::Round($d3.ToOADate())

    $cell.Value2 = $newSerial
}
